$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Non OCRED - OCR Problems")

# Add a new row of data (row 29) below the existing table
# Values are written in the same order the shared-string table records them
# (A29, then C29, then B29) so the underlying xl/sharedStrings.xml matches.
$ws.Range("A29").Value = "9,16,26"
$ws.Range("C29").Value = "Full article, contrasts"
$ws.Range("B29").Value = "https://demo.humlab.umu.se/courier/074055engo.pdf"

# Match the cursor/selection position recorded in the saved file
$ws.Activate()
$ws.Range("H25").Select()
